$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 corresponds to Johnathan Uselmann (Count=6). Fill in his
# previously-blank Favorite Ice Cream / Favorite Pizza Toppings cells.
$ws.Range("E9").Value = "Chocolate Reces' Penut butter"
$ws.Range("F9").Value = "Black Olives"
